$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Hyperlinks.Delete()
$r = $ws.Range("A1")
$r.Hyperlinks.Add($r, "https://first.example.com")
Write-Host "done"
